# Add a new worksheet "currency_movements" after the last existing sheet
# (mirrors: <sheet name="currency_movements" sheetId="10" r:id="rId7"/>)
$wb = $excel.ActiveWorkbook

$lastIndex = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($lastIndex)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "currency_movements"

# ---- Header row -------------------------------------------------------
$ws.Cells.Item(1, 1).Value = "date"
$ws.Cells.Item(1, 2).Value = "buy_date"
$ws.Cells.Item(1, 3).Value = "amount"
$ws.Cells.Item(1, 4).Value = "currency"
$ws.Cells.Item(1, 5).Value = "comment"
$headerRange = $ws.Range("A1:E1")
$headerRange.Font.Bold = $true
$headerRange.Font.Name = "Calibri"
$headerRange.Font.Size = 12

# ---- Data rows ----------------------------------------------------------
# Row 2 - USD, 100, first movement (fees excluded -> blank comment cell)
$ws.Cells.Item(2, 1).Value = 45292
$ws.Cells.Item(2, 2).Value = 44842
$ws.Cells.Item(2, 3).Value = 100
$ws.Cells.Item(2, 4).Value = "USD"
# keep E2 present (empty) but formatted like its row neighbours
$ws.Cells.Item(2, 3).Font.Name = "Calibri"
$ws.Cells.Item(2, 3).Font.Size = 12
$ws.Cells.Item(2, 4).Font.Name = "Calibri"
$ws.Cells.Item(2, 4).Font.Size = 12
$ws.Cells.Item(2, 5).Font.Name = "Calibri"
$ws.Cells.Item(2, 5).Font.Size = 12

# Row 3
$ws.Cells.Item(3, 1).Value = 45293
$ws.Cells.Item(3, 2).Value = 45293
$ws.Cells.Item(3, 3).Value = 100
$ws.Cells.Item(3, 4).Value = "EUR"

# Row 4 (manually corrected entry -> distinct font colour)
$ws.Cells.Item(4, 1).Value = 45294
$ws.Cells.Item(4, 2).Value = 45294
$ws.Cells.Item(4, 3).Value = -100
$ws.Cells.Item(4, 4).Value = "EUR"

# Row 5 (manually corrected entry -> distinct font colour)
$ws.Cells.Item(5, 1).Value = 45294
$ws.Cells.Item(5, 2).Value = 45294
$ws.Cells.Item(5, 3).Value = -100
$ws.Cells.Item(5, 4).Value = "USD"

# ---- Number formats -----------------------------------------------------
# date / buy_date columns -> same custom date format used elsewhere in the
# workbook (escaped dashes so the existing numFmtId 164 is reused).
$ws.Range("A2:B5").NumberFormat = "yyyy\-mm\-dd;@"

# Newly introduced font (black RGB instead of theme color) for the last two
# (manually corrected) rows.
$lastRowsDates = $ws.Range("A4:B5")
$lastRowsDates.Font.Color = 0
$lastRowsDates.Font.Name = "Calibri"
$lastRowsDates.Font.Size = 12

# ---- Move the tab-selected / active-sheet state ------------------------
# Before the edit, "currency_conversions" was the active / selected tab.
# After the edit it keeps a plain selection (no longer the tab in focus)
# and the new sheet becomes the active tab.
$prevActive = $wb.Worksheets.Item("currency_conversions")
$prevActive.Range("E33").Select() | Out-Null

$ws.Activate()
